# "Wii e Xbox perdidos" - add missing Xbox (Windows Store) and Wii (Mídia Física) games.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- XBOX / Windows Store block: rows 785-789 ---

# 1) Names for the four Gears of War entries (typed as a column first).
$ws.Cells.Item(785, 1).Value = "Gears of War"
$ws.Cells.Item(786, 1).Value = "Gears of War 2"
$ws.Cells.Item(787, 1).Value = "Gears of War 3"
$ws.Cells.Item(788, 1).Value = "Gears of War Judgement"

# 2) Logo URLs for those same four entries (typed as a column afterwards).
$ws.Cells.Item(785, 2).Value = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcTuCOqWFQ10b3zOq9O8OyfQtS2LD_sd0rP3Ip4iRX4oep69Xwc0"
$ws.Cells.Item(786, 2).Value = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcSVNIHJw8-XsYsAGaNK6tNtE0UAD-Qrz1Cp7SfKbqnK3ichiXTC"
$ws.Cells.Item(787, 2).Value = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcTWb0ricljrhot41Tm1hBKid157iphAndcjfeSnPBGDwY7Ks-g7Fw"
$ws.Cells.Item(788, 2).Value = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcRlvwhwLXvywBXo7kOkxYBOrTggBs-xs2FtYtzWGHhoVm7Jz3X55Q"

# 3) Assassins Creed Unity entered as a normal single row (name + url together).
$ws.Cells.Item(789, 1).Value = "Assassins Creed Unity"
$ws.Cells.Item(789, 2).Value = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcTcp6x7rM0N3sre7cHNSiv66VRKWGz_npaG7MbFMq06iloqK3-z"

# 4) Store / system / disabled columns for the whole XBOX block.
for ($r = 785; $r -le 789; $r++) {
    $ws.Cells.Item($r, 3).Value = "Windows Store"
    $ws.Cells.Item($r, 5).Value = "XBOX One"
}

# --- Wii / Mídia Física block: rows 790-802 ---

# 1) Names typed first, in entry order (rows 790-802, skipping 794 which already
#    held "resident evil 4" from elsewhere and is reused as-is).
$wiiNameEntryOrder = @(
    @{ Row = 790; Name = "Mario Kart" },
    @{ Row = 791; Name = "Donkey Kong Country Returns" },
    @{ Row = 792; Name = "The Legend of Zelda Twilight Princess" },
    @{ Row = 793; Name = "The Last Story" },
    @{ Row = 795; Name = "Wii Sports Resort" },
    @{ Row = 796; Name = "New Super Mario Bros Wii" },
    @{ Row = 797; Name = "Super Mario Galaxy" },
    @{ Row = 798; Name = "Pikmin 2" },
    @{ Row = 799; Name = "Super Paper Mario" },
    @{ Row = 800; Name = "The Legend of Zelda Skyward Sword" },
    @{ Row = 801; Name = "Wii Sports" },
    @{ Row = 802; Name = "Super Mario Galaxy 2" }
)
foreach ($d in $wiiNameEntryOrder) {
    $ws.Cells.Item($d.Row, 1).Value = $d.Name
}
$ws.Cells.Item(794, 1).Value = "resident evil 4 "

# 2) Logo URL / store / system / disabled typed per row, in row order 790-802.
$wiiRowData = @(
    @{ Row = 790; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcTJhjYxcs0DSUxaCjTRTRQXeJtPzjnUTbQe7rqckOBYUymXuur2" },
    @{ Row = 791; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcSbwqWLNZS7t5txSQ_J4QnOTTUQp19tdyf6lIsJwP6XRbmMUihV" },
    @{ Row = 792; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcR7kJCv5tK9CHK6WCEcpv7cFyRHeRJehX4it4-oFyVejYawpR-gjw" },
    @{ Row = 793; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcSMnHlTIyYkstE9M9XFTVz9HmTdwRgP0nHjj_9C-OGNkSuo16kTlg" },
    @{ Row = 794; Url = "http://media.steampowered.com/steamcommunity/public/images/apps/254700/532d72710af44f29cc123c5796e95e0382461ee5.jpg" },
    @{ Row = 795; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcQZAb6EFYyMel2Vs6HtrFksXLTAROFPeQ9tpS5iNYubKS1LImsb" },
    @{ Row = 796; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcSIHeicEMJFtogkX0mzINLEi2P_KJ3B4JtXxh44-UcXe-gbCgJW" },
    @{ Row = 797; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcTFHG5N3acPGb1NW2R3LYG-nY5tvZnDLEs2hYsrxS-atIsS0UmLdQ" },
    @{ Row = 798; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcTTTwDR4q2bqlSOuP3q3M8eRVJXFH6FKzfsQH6yX1d7JWttSbDD" },
    @{ Row = 799; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcQLpB6nukw_R5S_f19hKMjs1Jiadp74wSofWxvDszN9V60T2S-C" },
    @{ Row = 800; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcRBOH61IB-_qkPZI7sVkVi_3yKdpG-KL-F3qUvVOjIlI21Zl8WqNQ" },
    @{ Row = 801; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcS7v5-aEv5YciskPBhO__KnXNelmyFdr-4AyqMYy63A-DwazsKAWA" },
    @{ Row = 802; Url = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcTdpmEhrBIB59nSV1jD6wqI1jUzwwF12JSMiEWt2g_FXCUOoXxs" }
)
foreach ($d in $wiiRowData) {
    $ws.Cells.Item($d.Row, 2).Value = $d.Url
    $ws.Cells.Item($d.Row, 3).Value = "Mídia Física"
    $ws.Cells.Item($d.Row, 5).Value = "Wii"
}

# 3) Sort just the name column for the Wii block (A790:A802), like the original edit did.
$rng = $ws.Range("A790:A802")
$rng.Sort($ws.Range("A790"))

# --- "disabled" column (F) for every new row: must stay literal text "false" ---
# (A plain Value assignment of "false" gets typed as a boolean by this engine, so
# copy/paste-special a cell that's already a text "false" to keep it a string.)
$ws.Range("F2").Copy()
for ($r = 785; $r -le 802; $r++) {
    $ws.Cells.Item($r, 6).PasteSpecial(-4163)
}

# --- Update the current selection to match where the user ended up. ---
$ws.Range("B803").Select()
